$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.469.26"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.092.96"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'330.00"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.5217"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("D8").Value = "'0.4417"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'54.09"
$ws.Range("E9").Value = "  +17.49%  "
$ws.Range("D10").Value = "'0.08933"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("D12").Value = "'24.23"
$ws.Range("E12").Value = "  -4.52%  "
$ws.Range("D13").Value = "2.101.02"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "'6.679"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'7.700"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "'95.90"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "'0.00001122"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'0.06610"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'19.10"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'6.258"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "30.505.90"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").Value = "'12.32"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "2.341.06"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "'22.21"
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").Value = "'2.569"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "'163.82"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'131.55"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'1.190"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "'0.1069"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'1.646"
$ws.Range("E33").Value = "  +6.97%  "
$ws.Range("D34").Value = "'6.150"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "'3.901"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "'10.13"
$ws.Range("E36").Value = "  +5.47%  "
$ws.Range("D37").Value = "'0.02556"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "'0.06804"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'5.467"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "'12.61"
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").Value = "'0.2251"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").Value = "'0.6866"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'1.250"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'13.97"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "'0.6322"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'2.194"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").Value = "'3.630"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("D51").Value = "'81.59"
$ws.Range("E51").Value = "  -2.19%  "
